$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text-safe cell updates (won't be misinterpreted as numbers) ---
$ws.Range('D2').Value = '30.336.82'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').Value = '1.931.56'
$ws.Range('E3').Value = '  +0.23%  '
$ws.Range('E5').Value = '  +2.14%  '
$ws.Range('E6').Value = '  +0.46%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  +0.82%  '
$ws.Range('E9').Value = '  +4.69%  '
$ws.Range('E10').Value = '  +5.48%  '
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('E12').Value = '  +1.83%  '
$ws.Range('D13').Value = '1.929.57'
$ws.Range('E13').Value = '  +0.14%  '
$ws.Range('E14').Value = '  +0.63%  '
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('E16').Value = '  +2.15%  '
$ws.Range('D17').Value = '30.313.52'
$ws.Range('E17').Value = '  +0.21%  '
$ws.Range('E18').Value = '  -3.05%  '
$ws.Range('E19').Value = '  +2.57%  '
$ws.Range('E20').Value = '  -0.55%  '
$ws.Range('D21').Value = '2.184.03'
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('E24').Value = '  +1.06%  '
$ws.Range('E25').Value = '  +0.91%  '
$ws.Range('E26').Value = '  +3.20%  '
$ws.Range('E27').Value = '  +1.59%  '
$ws.Range('E28').Value = '  +1.74%  '
$ws.Range('E29').Value = '  -3.41%  '
$ws.Range('E30').Value = '  -0.12%  '
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('E32').Value = '  +0.38%  '
$ws.Range('E33').Value = '  +0.28%  '
$ws.Range('E34').Value = '  +3.04%  '
$ws.Range('E35').Value = '  +6.66%  '
$ws.Range('E36').Value = '  +1.24%  '
$ws.Range('E37').Value = '  +0.93%  '
$ws.Range('E38').Value = '  +1.40%  '
$ws.Range('E39').Value = '  -0.35%  '
$ws.Range('E40').Value = '  -0.89%  '
$ws.Range('E41').Value = '  -1.94%  '
$ws.Range('E42').Value = '  +1.82%  '
$ws.Range('E43').Value = '  +1.63%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('E44').Value = '  +1.12%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('E46').Value = '  -0.94%  '
$ws.Range('E47').Value = '  +0.58%  '
$ws.Range('E48').Value = '  +2.15%  '
$ws.Range('E49').Value = '  +1.10%  '
$ws.Range('E50').Value = '  +2.82%  '
$ws.Range('E51').Value = '  +1.90%  '

# --- Numeric-looking strings that must stay text: use NumberFormat trick ---
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '251.38'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.7185'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3267'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '27.54'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7985'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08083'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.64'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.78'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '252.03'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008132'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.784'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.919'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.718'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.16'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.319'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1280'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.359'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.544'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.431'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.197'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05212'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.269'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7478'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.765'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.799'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '79.06'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.434'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4523'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.029'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8411'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.001'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '101.77'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.760'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.407'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06078'
$ws.Range('D50').Style = 'Normal'
